# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces the target commit:
#   1. Slide 16's table (graphicFrame) gets a new table style
#      (tableStyleId {3B1D688B-650D-4EE3-8E1F-6868BAD2C909} ->
#       {C9E6C852-2527-42AD-A317-34126E4201F0}).
#   2. The presentation's theme (ppt/theme/theme1.xml, the theme used by
#      the slide master / "Integral" design) is recolored to the stock
#      PowerPoint "Office" palette - i.e. the 12-slot theme color scheme
#      that used to live only on the notes master's theme is now also the
#      slide master's theme color scheme (dk2/lt2/accent1-6/hlink/folHlink
#      all change; dk1/lt1 stay black/white).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: turn an "RRGGBB" hex string into the BGR-packed integer that
# the PowerPoint object model's ColorFormat/ThemeColor `.RGB` property
# uses (same convention as VBA's RGB() function).
# ---------------------------------------------------------------------
function HexToComRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# ---------------------------------------------------------------------
# 1) Table style change on slide 16's table.
# ---------------------------------------------------------------------
$targetSlideIndex = 16
$oldTableStyleId = "{3B1D688B-650D-4EE3-8E1F-6868BAD2C909}"
$newTableStyleId = "{C9E6C852-2527-42AD-A317-34126E4201F0}"

$slide = $p.Slides.Item($targetSlideIndex)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $table = $shape.Table
        if ($table.Style -eq $oldTableStyleId) {
            $table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Recolor the slide master's theme color scheme to the stock
#    "Office" palette (dk1/lt1 are unchanged: black/white).
# ---------------------------------------------------------------------
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$master = $p.SlideMaster
$themeColorScheme = $master.Theme.ThemeColorScheme
for ($idx = 1; $idx -le $officeThemeColors.Count; $idx++) {
    $themeColorScheme.Item($idx).RGB = HexToComRgb $officeThemeColors[$idx - 1]
}
